$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HI")

# Fix typo in the Roman transliteration name: "Achuytam Keshavam" -> "Achyutam Keshavam"
$ws.Range("C2").Value = "Achyutam Keshavam"

# Update the active selection on the sheet to C3
$ws.Activate()
$ws.Range("C3").Select()
